$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for [$findText]"
    }
    return $ok
}

# --- Title page: version + date ---
ReplaceText "Version1.2" "Version1.3" | Out-Null
ReplaceText "December 20, 2013" "January 10, 2014" | Out-Null

# --- SPUID paragraph text corrections ---
ReplaceText "submitter_id" "spuid_namespace" | Out-Null
ReplaceText " attribute that specifies submitter. The values of " " attribute that is unique for each submitter. The values of " | Out-Null
ReplaceText " are rom controlled vocabulary" " are from controlled vocabulary" | Out-Null

# --- Move the _GoBack bookmark into the SPUID paragraph, right after "SPUID i" ---
$findRange = $d.Content
$found = $findRange.Find.Execute("SPUID is used to link")
if ($found) {
    $insertPos = $findRange.Start + 7
    $bmRange = $d.Range($insertPos, $insertPos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# --- Footer: cached PAGE field result 1 -> 4 ---
$footer = $d.Sections(1).Footers(1)
$fok = $footer.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2)
Write-Output "Footer replace: $fok"

Write-Output "DONE"
